$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 4 (old rows 4-5 shift down to 7-8)
$ws.Rows("4:6").Insert()

# New data for the three inserted rows (Brazil - Serie A Betano matches)
$newRows = @(
    @{
        A="C6pXYSIf"; B="20/11/2024"; C="16:30"; D="BRAZIL - SERIE A BETANO"; E="Athletico-PR"; F="Atletico GO";
        G=1.62; H=3.7; I=5.75; J=2.25; K=2.2; L=5.5; M=1.06; N=10; O=1.3; P=3.5; Q=2.02; R=1.88; S=1.4; T=2.75; U=1.95; V=1.8;
        W=6.5; X=7.5; Y=8.5; Z=12; AA=13; AB=29; AC=9.5; AD=7; AE=17; AF=51; AG=351; AH=13; AI=29; AJ=17; AK=51; AL=41; AM=41;
        AN=3.5; AO=8.5; AP=21; AQ=29; AR=51; AS=151; AT=2.75; AU=8.5; AV=51; AW=7; AX=29; AY=34; AZ=101; BA=126; BB=301; BC=126; BD=126
    },
    @{
        A="IqoTZ83l"; B="20/11/2024"; C="16:30"; D="BRAZIL - SERIE A BETANO"; E="Bragantino"; F="Sao Paulo";
        G=2.63; H=3; I=3; J=3.4; K=1.95; L=3.75; M=1.1; N=7; O=1.44; P=2.75; Q=2.4; R=1.53; S=1.53; T=2.38; U=2; V=1.75;
        W=7; X=11; Y=11; Z=26; AA=23; AB=41; AC=7; AD=5.5; AE=17; AF=67; AG=451; AH=7.5; AI=13; AJ=12; AK=29; AL=29; AM=41;
        AN=4.5; AO=15; AP=29; AQ=51; AR=81; AS=251; AT=2.38; AU=9; AV=67; AW=4.75; AX=17; AY=29; AZ=51; BA=101; BB=301; BC=126; BD=126
    },
    @{
        A="8YwtX6m7"; B="20/11/2024"; C="16:30"; D="BRAZIL - SERIE A BETANO"; E="Criciuma"; F="Vitoria";
        G=2.3; H=3.2; I=3.25; J=3; K=2; L=4; M=1.08; N=8; O=1.4; P=3; Q=2.2; R=1.65; S=1.5; T=2.5; U=1.95; V=1.8;
        W=7; X=10; Y=9.5; Z=21; AA=21; AB=34; AC=8; AD=6; AE=15; AF=51; AG=351; AH=8.5; AI=15; AJ=12; AK=34; AL=29; AM=41;
        AN=4.33; AO=13; AP=26; AQ=41; AR=67; AS=201; AT=2.5; AU=8.5; AV=67; AW=5; AX=19; AY=29; AZ=67; BA=101; BB=251; BC=126; BD=126
    }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD")

$rowIdx = 4
foreach ($rowData in $newRows) {
    foreach ($col in $cols) {
        $ws.Range("$col$rowIdx").Value = $rowData[$col]
    }
    $rowIdx++
}
